# Actualización 11 de Mayo - Mañana
# Populate the "Rescatables" sheet with 4 rescatable-student rows (rows 2-5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Column A - NC (control numbers)
$ws.Range("A2").Value = 20330051920202
$ws.Range("A3").Value = 18330051920237
$ws.Range("A4").Value = 18330051920248
$ws.Range("A5").Value = 18330051920245

# Column B - Paterno
$ws.Range("B2").Value = "DE JESUS"
$ws.Range("B3").Value = "CANCINO"
$ws.Range("B4").Value = "HERNANDEZ"
$ws.Range("B5").Value = "DE JESUS"

# Column C - Materno
$ws.Range("C2").Value = "CASTILLO"
$ws.Range("C3").Value = "GUERRA"
$ws.Range("C4").Value = "DOLORES"
$ws.Range("C5").Value = "ISIDRO"

# Column D - Nombres
$ws.Range("D2").Value = "ITZEL"
$ws.Range("D3").Value = "DANIEL"
$ws.Range("D4").Value = "JOEL EDUARDO"
$ws.Range("D5").Value = "MONSERRAT"

# Column E - Nombre_Largo (materia)
$ws.Range("E2").Value = "TOMA MUESTRAS BIOLÓGICAS"
$ws.Range("E3").Value = "ANALIZA SANGRE CON BASE EN TÉCNICAS DE QUÍMICA CLÍNICA"
$ws.Range("E4").Value = "ANALIZA SANGRE CON BASE EN TÉCNICAS DE QUÍMICA CLÍNICA"
$ws.Range("E5").Value = "ANALIZA SANGRE CON BASE EN TÉCNICAS DE QUÍMICA CLÍNICA"

# Column F - Grupo
$ws.Range("F2").Value = "2ALCV"
$ws.Range("F3").Value = "6ALCV"
$ws.Range("F4").Value = "6ALCV"
$ws.Range("F5").Value = "6ALCV"

# Column G - Reprobadas
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2
